{"js": "// Insert the missing contact-info paragraph right after the name\n// heading (\"Dheeraj Chand\") and before the \"PROFESSIONAL SUMMARY\"\n// heading, matching the target diff.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// Locate the first paragraph (the centered name heading \"Dheeraj Chand\").\nconst nameParagraph = paragraphs.items[0];\n\nconst contactText =\n  \"202.550.7110 | dheeraj.chand@gmail.com | https://www.dheerajchand.com | https://www.linkedin.com/in/dheerajchand/ | Austin, TX\";\n\n// Insert a brand-new paragraph directly after the name heading. The new\n// paragraph initially inherits the name heading's bold/28pt run\n// formatting, so clear() it (drops the inherited direct formatting) and\n// then insert the plain contact text.\nconst contactParagraph = nameParagraph.insertParagraph(\"\", Word.InsertLocation.after);\ncontactParagraph.alignment = Word.Alignment.centered;\ncontactParagraph.clear();\nawait context.sync();\n\ncontactParagraph.insertText(contactText, Word.InsertLocation.start);\n\nawait context.sync();\n", "ps1": "# Fix contact information missing from short resumes: add the contact\n# info line as its own centered paragraph right after the \"Dheeraj\n# Chand\" name heading and before the \"PROFESSIONAL SUMMARY\" heading.\n\n$d = $word.ActiveDocument\n\n$contactText = \"202.550.7110 | dheeraj.chand@gmail.com | https://www.dheerajchand.com | https://www.linkedin.com/in/dheerajchand/ | Austin, TX\"\n\n$rng = $d.Content\n$rng.Find.ClearFormatting()\n$rng.Find.Text = \"Dheeraj Chand\"\n$rng.Find.Replacement.ClearFormatting()\n# \"^p\" is Word's paragraph-mark code, so the replacement splits the name\n# heading into its own paragraph followed by a brand-new paragraph that\n# holds the contact info (inheriting the centered alignment, but none of\n# the name's bold/28pt direct character formatting).\n$rng.Find.Replacement.Text = \"Dheeraj Chand^p\" + $contactText\n$rng.Find.Execute($null, $null, $null, $null, $null, $null, $null, $null, $null, $null, 2)\n"}
